# updated legacy GSC export data
#  - Append 4 new daily rows (2025-11-14 .. 2025-11-17) to the "Chart" sheet.
#  - Update the "Videos" count on the "Table" sheet summary row (9 -> 8).

$wb = $excel.ActiveWorkbook

# --- Chart sheet: append new rows 88-91 ------------------------------------
$chart = $wb.Worksheets.Item("Chart")

$newRows = @(
    @{ Row = 88; Date = "2025-11-14"; NoVideo = 8.0; VideoIndexed = 0.0; Impressions = 0.0 },
    @{ Row = 89; Date = "2025-11-15"; NoVideo = 8.0; VideoIndexed = 0.0; Impressions = 0.0 },
    @{ Row = 90; Date = "2025-11-16"; NoVideo = 8.0; VideoIndexed = 0.0; Impressions = 0.0 },
    @{ Row = 91; Date = "2025-11-17"; NoVideo = 8.0; VideoIndexed = 0.0; Impressions = $null }
)

foreach ($r in $newRows) {
    $dateCell = $chart.Cells.Item($r.Row, 1)
    # leading apostrophe forces text entry so the ISO date string isn't
    # auto-converted into a date serial number
    $dateCell.Value = "'" + $r.Date
    $dateCell.ClearFormats()

    $chart.Cells.Item($r.Row, 2).Value = $r.NoVideo
    $chart.Cells.Item($r.Row, 3).Value = $r.VideoIndexed

    $impCell = $chart.Cells.Item($r.Row, 4)
    if ($r.Impressions -eq $null) {
        # last day has no impressions figure recorded yet -> blank text cell
        $impCell.Value = "'"
        $impCell.ClearFormats()
    } else {
        $impCell.Value = $r.Impressions
    }
}

# --- Table sheet: update Videos count --------------------------------------
$table = $wb.Worksheets.Item("Table")
$table.Range("C2").Value = 8.0
